$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# New "Food" entries logged on 2019-04-01 (activity log continues past row 153)
$newRows = @(
    @{ Row = 154; Start = 43556.311805555553; Comment = "Whole Wheat Bread" },
    @{ Row = 155; Start = 43556.506944444445; Comment = "Beans + rice" },
    @{ Row = 156; Start = 43556.5625;          Comment = "Kombucha" },
    @{ Row = 157; Start = 43556.645833333336; Comment = "Blueberry Kefir" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $prev = $r - 1

    # Inherit the row-above formatting (date style on Start, calculated-column
    # style on Z) before filling in this row's own values/formula.
    $ws.Range("A$($prev):E$($prev)").Copy()
    $ws.Range("A$($r):E$($r)").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($r, 2).Clear()                   # no "End" value for Food rows

    $ws.Cells.Item($r, 1).Value = $item.Start
    $ws.Cells.Item($r, 3).Value = "Food"
    $ws.Cells.Item($r, 4).Value = $item.Comment
    $ws.Cells.Item($r, 5).Formula = '=IF(Table2[[#This Row],[Activity]]="Sleep",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,"NA")'
}

# Grow the table to cover the newly-entered rows
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:E157"))

# Scroll/selection state left by the editor after entering the last row
$ws.Application.ActiveWindow.ScrollRow = 127
$ws.Range("A158").Select()
